# Adds the submitter's name, email and repo link as a second row, turning
# the email and repo link into real hyperlinks (with the matching
# "Hyperlink" cell style), resizes the columns to fit the new content, and
# moves the active selection - mirroring the author's own edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data -----------------------------------------------------------
$ws.Range("A2").Value = "mai medhat mohamed qandil"

# Hyperlinks.Add both writes the cell text/display value and wires up the
# relationship + "Hyperlink" style (underline, theme color 10) in one go.
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:mkandel2025@gmail.com", "", "", "mkandel2025@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/maiqandil/Security-Task.git", "", "", "https://github.com/maiqandil/Security-Task.git")

# --- Column widths ----------------------------------------------------------
# Widened to fit the newly-added name / email / repo-link text.
$ws.Columns.Item(1).ColumnWidth = 26.5
$ws.Columns.Item(2).ColumnWidth = 35
$ws.Columns.Item(3).ColumnWidth = 53

# --- Selection --------------------------------------------------------------
$ws.Range("C7").Select() | Out-Null
